$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column headers (A-T) stay the same; only rows 2-9 of data change.
$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

$data = @(
    @("FAPs", "Artn", "Ret", "ECs", 3, 1, 1.415594, 4.246782, 0.9753967820130081, 0.9753967820130081, 3, 1, 3.632885, 10.898655, 0.4372849566404539, 0.4372849566404539, 5.142690208689999, 46.28421187820999, 0.4265263395297965, 0.4265263395297965),
    @("FAPs", "Artn", "Ret", "FAPs", 3, 1, 1.415594, 4.246782, 0.9753967820130081, 0.9753967820130081, 3, 1, 3.884996, 11.654988, 0.4676311822169809, 0.4676311822169809, 5.499577027623999, 49.496193248616, 0.4561259503033818, 0.4561259503033818),
    @("FAPs", "Artn", "Ret", "M1", 3, 1, 1.415594, 4.246782, 0.9753967820130081, 0.9753967820130081, 1, 0.3333333333333333, 0.004340666666666667, 0.013022, 0.0005224795816889323, 0.0005224795816889323, 0.006144621689333333, 0.055301595204, 0.0005096249026468871, 0.0005096249026468871),
    @("FAPs", "Artn", "Ret", "sCs", 3, 1, 1.415594, 4.246782, 0.9753967820130081, 0.9753967820130081, 3, 1, 0.7855989999999999, 2.356797, 0.09456138156087625, 0.09456138156087625, 1.112089230806, 10.008803077254, 0.09223486727718289, 0.09223486727718289),
    @("sCs", "Artn", "Ret", "ECs", 1, 0.3333333333333333, 0.03570666666666667, 0.10712, 0.0246032179869919, 0.0246032179869919, 3, 1, 3.632885, 10.898655, 0.4372849566404539, 0.4372849566404539, 0.1297182137333333, 1.1674639236, 0.01075861711065739, 0.01075861711065739),
    @("sCs", "Artn", "Ret", "FAPs", 1, 0.3333333333333333, 0.03570666666666667, 0.10712, 0.0246032179869919, 0.0246032179869919, 3, 1, 3.884996, 11.654988, 0.4676311822169809, 0.4676311822169809, 0.1387202571733333, 1.24848231456, 0.01150523191359911, 0.01150523191359911),
    @("sCs", "Artn", "Ret", "M1", 1, 0.3333333333333333, 0.03570666666666667, 0.10712, 0.0246032179869919, 0.0246032179869919, 1, 0.3333333333333333, 0.004340666666666667, 0.013022, 0.0005224795816889323, 0.0005224795816889323, 0.0001549907377777778, 0.00139491664, 0.00001285467904204515, 0.00001285467904204515),
    @("sCs", "Artn", "Ret", "sCs", 1, 0.3333333333333333, 0.03570666666666667, 0.10712, 0.0246032179869919, 0.0246032179869919, 3, 1, 0.7855989999999999, 2.356797, 0.09456138156087625, 0.09456138156087625, 0.02805112162666667, 0.25246009464, 0.002326514283693355, 0.002326514283693355)
)

$rowIndex = 2
foreach ($row in $data) {
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $addr = "$($cols[$i])$rowIndex"
        $ws.Range($addr).Value = $row[$i]
    }
    $rowIndex++
}
